$wb = $excel.ActiveWorkbook

# --- ALC ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H11").Value = 335.76923
$ws.Range("I11").Value = 335.76923
$ws.Range("K11").Value = 335.76923
$ws.Range("M11").Value = -195.76923
$ws.Range("H88").Value = 1897.8334
$ws.Range("N88").Value = -2769.4
$ws.Range("J88").Value = 1957.4
$ws.Range("L88").Value = 1957.4
$ws.Range("K88").Value = 1600
$ws.Range("M88").Value = -1194
$ws.Range("I88").Value = 1600
$ws.Range("I91").Value = 1600
$ws.Range("N91").Value = -4765.4
$ws.Range("L91").Value = 1957.4
$ws.Range("M91").Value = -196
$ws.Range("K91").Value = 1600
$ws.Range("J91").Value = 1957.4
$ws.Range("H91").Value = 1897.8334
$ws.Range("H93").Value = 30000
$ws.Range("J93").Value = 30000
$ws.Range("N93").Value = -34992
$ws.Range("L93").Value = 30000
$ws.Range("N106").Value = -3859.6667
$ws.Range("H106").Value = 31252024
$ws.Range("J106").Value = 2597.6667
$ws.Range("M106").Value = -38462799
$ws.Range("K106").Value = 38463430
$ws.Range("I106").Value = 38463430
$ws.Range("L106").Value = 2597.6667
$ws.Range("L132").Value = 10699.3329
$ws.Range("I132").Value = 1475.76
$ws.Range("K132").Value = 4427.28
$ws.Range("J132").Value = 3566.4443
$ws.Range("M132").Value = -1897.28
$ws.Range("N132").Value = -15759.3329
$ws.Range("H132").Value = 2029.1765

# --- ARM ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("N61").Value = -5213.5
$ws.Range("M61").Value = -2426.2974
$ws.Range("L61").Value = 4789.5
$ws.Range("I61").Value = 2638.2974
$ws.Range("K61").Value = 2638.2974
$ws.Range("H61").Value = 2938.465
$ws.Range("J61").Value = 4789.5
$ws.Range("I74").Value = 2057.7334
$ws.Range("H74").Value = 6645.3076
$ws.Range("K74").Value = 2057.7334
$ws.Range("M74").Value = -1183.7334
$ws.Range("I77").Value = 2057.7334
$ws.Range("H77").Value = 6645.3076
$ws.Range("K77").Value = 10288.667
$ws.Range("M77").Value = -5920.667000000001
$ws.Range("I122").Value = 2440.9375
$ws.Range("H122").Value = 2552.75
$ws.Range("K122").Value = 7322.8125
$ws.Range("M122").Value = -4872.8125
$ws.Range("L136").Value = 14368.5
$ws.Range("H136").Value = 2938.465
$ws.Range("N136").Value = -19468.5
$ws.Range("I136").Value = 2638.2974
$ws.Range("J136").Value = 4789.5
$ws.Range("M136").Value = -5364.8922
$ws.Range("K136").Value = 7914.8922

# --- BSM ---
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("K94").Value = 1749.4286
$ws.Range("H94").Value = 1935.7241
$ws.Range("M94").Value = -1298.4286
$ws.Range("J94").Value = 2424.75
$ws.Range("N94").Value = -3326.75
$ws.Range("I94").Value = 1749.4286
$ws.Range("L94").Value = 2424.75

# --- CRP ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("M22").Value = -440
$ws.Range("L22").Value = 0
$ws.Range("K22").Value = 790
$ws.Range("I22").Value = 790
$ws.Range("J22").Value = 0
$ws.Range("N22").ClearContents()
$ws.Range("H22").Value = 790
$ws.Range("K58").Value = 2312.8206
$ws.Range("L58").Value = 3170.5557
$ws.Range("I58").Value = 2312.8206
$ws.Range("J58").Value = 3170.5557
$ws.Range("N58").Value = -3576.5557
$ws.Range("M58").Value = -2109.8206
$ws.Range("H58").Value = 2473.6458
$ws.Range("H99").Value = 3698.111
$ws.Range("K99").Value = 3946.6
$ws.Range("J99").Value = 3387.5
$ws.Range("I99").Value = 3946.6
$ws.Range("L99").Value = 3387.5
$ws.Range("M99").Value = -2448.6
$ws.Range("N99").Value = -6383.5
$ws.Range("H122").Value = 957.5
$ws.Range("J122").Value = 1098.25
$ws.Range("L122").Value = 3294.75
$ws.Range("N122").Value = -8194.75
$ws.Range("I126").Value = 3946.6
$ws.Range("L126").Value = 10162.5
$ws.Range("M126").Value = -9369.799999999999
$ws.Range("J126").Value = 3387.5
$ws.Range("H126").Value = 3698.111
$ws.Range("K126").Value = 11839.8
$ws.Range("N126").Value = -15102.5
$ws.Range("H133").Value = 25375.334
$ws.Range("J133").Value = 25375.334
$ws.Range("N133").Value = -30435.334
$ws.Range("L133").Value = 25375.334
$ws.Range("M134").Value = -40992.894
$ws.Range("N134").Value = -140855.331
$ws.Range("J134").Value = 45261.777
$ws.Range("K134").Value = 43527.894
$ws.Range("L134").Value = 135785.331
$ws.Range("I134").Value = 14509.298
$ws.Range("H134").Value = 19451.66
$ws.Range("L136").Value = 9511.667099999999
$ws.Range("H136").Value = 2473.6458
$ws.Range("N136").Value = -14611.6671
$ws.Range("I136").Value = 2312.8206
$ws.Range("J136").Value = 3170.5557
$ws.Range("M136").Value = -4388.4618
$ws.Range("K136").Value = 6938.4618

# --- CUL ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("L2").Value = 343.09092
$ws.Range("M2").Value = -214793.29
$ws.Range("H2").Value = 20083.08
$ws.Range("N2").Value = -569.09092
$ws.Range("I2").Value = 35817.715
$ws.Range("J2").Value = 57.18182
$ws.Range("K2").Value = 214906.29
$ws.Range("J4").Value = 1000000
$ws.Range("K4").Value = 87375618
$ws.Range("M4").Value = -87375506
$ws.Range("L4").Value = 3000000
$ws.Range("N4").Value = -3000224
$ws.Range("I4").Value = 29125206
$ws.Range("H4").Value = 28455558
$ws.Range("I34").Value = 395
$ws.Range("H34").Value = 4631549
$ws.Range("L34").Value = 18525799.5
$ws.Range("K34").Value = 1185
$ws.Range("M34").Value = -1101
$ws.Range("N34").Value = -18525967.5
$ws.Range("J34").Value = 6175266.5
$ws.Range("H38").Value = 46.46154
$ws.Range("J38").Value = 44
$ws.Range("L38").Value = 132
$ws.Range("N38").Value = -826
$ws.Range("H39").Value = 1834.5
$ws.Range("J39").Value = 2001.4
$ws.Range("L39").Value = 6004.200000000001
$ws.Range("N39").Value = -6592.200000000001
$ws.Range("H55").Value = 2138
$ws.Range("J55").Value = 3325
$ws.Range("L55").Value = 9975
$ws.Range("N55").Value = -10329
$ws.Range("I92").Value = 1343.1428
$ws.Range("H92").Value = 1435.2
$ws.Range("K92").Value = 4029.4284
$ws.Range("M92").Value = -2781.4284
$ws.Range("I131").Value = 201018
$ws.Range("K131").Value = 603054
$ws.Range("H131").Value = 28199.744
$ws.Range("M131").Value = -598014
$ws.Range("L131").Value = 8355.882599999999
$ws.Range("N131").Value = -18435.8826
$ws.Range("J131").Value = 2785.2942

# --- GSM ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H80").Value = 5465.1665
$ws.Range("I80").Value = 2427.7144
$ws.Range("K80").Value = 2427.7144
$ws.Range("M80").Value = -1429.7144
$ws.Range("H83").Value = 5465.1665
$ws.Range("I83").Value = 2427.7144
$ws.Range("K83").Value = 12138.572
$ws.Range("M83").Value = -7146.572
$ws.Range("I95").Value = 0
$ws.Range("N95").Value = -33125.334
$ws.Range("J95").Value = 27633.334
$ws.Range("M95").ClearContents()
$ws.Range("K95").Value = 0
$ws.Range("H95").Value = 27633.334
$ws.Range("L95").Value = 27633.334
$ws.Range("H97").Value = 2684.8845
$ws.Range("J97").Value = 2908.5
$ws.Range("L97").Value = 2908.5
$ws.Range("I97").Value = 2585.5
$ws.Range("N97").Value = -3900.5
$ws.Range("M97").Value = -2089.5
$ws.Range("K97").Value = 2585.5
$ws.Range("I102").Value = 1330.2354
$ws.Range("H102").Value = 47620290
$ws.Range("K102").Value = 1330.2354
$ws.Range("M102").Value = 291.7646

# --- LTW ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H36").Value = 62500
$ws.Range("J36").Value = 62500
$ws.Range("N36").Value = -63624
$ws.Range("L36").Value = 62500
$ws.Range("I82").Value = 2352.3845
$ws.Range("N82").Value = -2809
$ws.Range("L82").Value = 2087
$ws.Range("M82").Value = -1991.3845
$ws.Range("K82").Value = 2352.3845
$ws.Range("H82").Value = 2219.6924
$ws.Range("J82").Value = 2087
$ws.Range("I85").Value = 2352.3845
$ws.Range("M85").Value = -1104.3845
$ws.Range("N85").Value = -4583
$ws.Range("H85").Value = 2219.6924
$ws.Range("K85").Value = 2352.3845
$ws.Range("L85").Value = 2087
$ws.Range("J85").Value = 2087
$ws.Range("H136").Value = 2848.861
$ws.Range("I136").Value = 2827.2144
$ws.Range("K136").Value = 8481.643199999999
$ws.Range("M136").Value = -5931.643199999999

# --- WVR ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("J81").Value = 5225
$ws.Range("M81").Value = -20813.834
$ws.Range("K81").Value = 21874.834
$ws.Range("I81").Value = 10937.417
$ws.Range("N81").Value = -12572
$ws.Range("H81").Value = 8081.2085
$ws.Range("L81").Value = 10450
$ws.Range("I84").Value = 10937.417
$ws.Range("H84").Value = 8081.2085
$ws.Range("K84").Value = 109374.17
$ws.Range("M84").Value = -104070.17
$ws.Range("J84").Value = 5225
$ws.Range("L84").Value = 52250
$ws.Range("N84").Value = -62858
$ws.Range("H95").Value = 93328.336
$ws.Range("J95").Value = 93328.336
$ws.Range("N95").Value = -98820.336
$ws.Range("L95").Value = 93328.336
$ws.Range("H97").Value = 35300
$ws.Range("J97").Value = 35300
$ws.Range("N97").Value = -37282
$ws.Range("L97").Value = 35300
$ws.Range("L132").Value = 4093.6665
$ws.Range("I132").Value = 1973.3043
$ws.Range("K132").Value = 5919.9129
$ws.Range("J132").Value = 1364.5555
$ws.Range("M132").Value = -3389.9129
$ws.Range("N132").Value = -9153.666499999999
$ws.Range("H132").Value = 1802.0938
